$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update standings table values (rows 6-9)
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = -5
$ws.Range("G6").Value = 2

$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 2

$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 2

$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = -5
$ws.Range("G9").Value = 2

# Move the selection to I10, matching the saved cursor position
$ws.Range("I10").Select()
